# Pinbeschreibung.xlsx update:
# Insert two new columns (Port: GPIOx / Pin: GPIO_PIN_x) after column B,
# pushing the former C..F columns to E..H. Fill the new columns for the
# rows that specify a GPIO port/pin (rows 8-12), and update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at C:D - this shifts old C,D,E,F -> E,F,G,H
$ws.Columns("C:D").Insert()

# New header cells for the inserted columns
$ws.Range("C1").Value = "Port: GPIOx"
$ws.Range("D1").Value = "Pin: GPIO_PIN_x"

# New GPIO port/pin data for rows 8-12
$ws.Range("C8").Value  = "GPIOA"
$ws.Range("D8").Value  = "GPIO_PIN_6"

$ws.Range("C9").Value  = "GPIOB"
$ws.Range("D9").Value  = "GPIO_PIN_5"

$ws.Range("C10").Value = "GPIOA"
$ws.Range("D10").Value = "GPIO_PIN_8"
$ws.Range("E10").Value = "Rotary encoder channel A"

$ws.Range("C11").Value = "GPIOA"
$ws.Range("D11").Value = "GPIO_PIN_9"
$ws.Range("E11").Value = "Rotary encoder channel B"

$ws.Range("C12").Value = "GPIOC"
$ws.Range("D12").Value = "GPIO_PIN_1"

# Match column width of the new columns to column B
$ws.Columns("C:D").ColumnWidth = 20.0

# Widen the "Use" column (now E) slightly to fit the new longer text
$ws.Columns("E").ColumnWidth = 22.4

# Update the selected cell
$ws.Range("C21").Select()
